# Daily update at 8 AM UTC
# Previously the last data row (row 19, date 2025-04-12) was styled with the
# "final row" date-only format. Now that a new day's data is appended, that
# older row reverts to the standard "YYYY-MM-DD HH:MM:SS" date format used by
# all the non-final rows, and the newly appended row 20 becomes the new
# "final row" with the date-only format.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 19 is no longer the last row - give it the regular date/time format
# used by every other (non-final) row in the table.
$ws.Range("A19").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Append the new day's row.
$ws.Range("A20").Value = 45760
$ws.Range("A20").NumberFormat = "YYYY-MM-DD"
$ws.Range("B20").Value = 79
$ws.Range("C20").Value = 79
$ws.Range("D20").Value = 76
